{"js": "// Collapse the split \"<id>...</id>\" runs (e.g. \"<id>\" + \"p098r_a1\" + \"</id>\")\n// into a single run reading \"<id>p098r_N</id>\" for each of the document's\n// three <id> placeholders, matching the formatting of the opening \"<id>\" run.\nconst replacements = [\n  { find: \"<id>p098r_a1</id>\", text: \"<id>p098r_1</id>\" },\n  { find: \"<id>p098r_a2</id>\", text: \"<id>p098r_2</id>\" },\n  { find: \"<id>p098r_a3</id>\", text: \"<id>p098r_3</id>\" }\n];\n\nfor (const { find, text } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(text, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Collapse the split \"<id>...</id>\" runs (e.g. \"<id>\" + \"p098r_a1\" + \"</id>\")\n# into a single run reading \"<id>p098r_N</id>\" for each of the document's\n# three <id> placeholders, matching the formatting of the opening \"<id>\" run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"<id>p098r_a1</id>\"; Replace = \"<id>p098r_1</id>\" },\n    @{ Find = \"<id>p098r_a2</id>\"; Replace = \"<id>p098r_2</id>\" },\n    @{ Find = \"<id>p098r_a3</id>\"; Replace = \"<id>p098r_3</id>\" }\n)\n\nforeach ($item in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $item.Find\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Execute($item.Find, $true, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2)\n}\n"}
